$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("components request")

# Set B2 to the new revision test string
$ws.Range("B2").Value = "test_rev_A2"

# Correct the RIGHT() formula's case for rev "0": extract 2 chars instead of 1
$ws.Range("C2").Formula = "=RIGHT(B2,IF(RIGHT(B2,1)=""0"",0,2))"

# Move the active selection to B3 to match the saved view state
$ws.Range("B3").Select() | Out-Null
